# Insert a new weekly price record at row 17 (shifting existing rows 17-130 down
# to 18-131) and populate it with the new entry's values.
#
# The new record mirrors the row directly above it (row 16: Feria Lagunitas de
# Puerto Montt / Albahaca / $/paquete / Region de Arica y Parinacota) but with a
# newer date (2022-07-07, Excel serial 44749) and its own volume figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..130 down to 18..131, creating a blank row 17.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new record's data.
$ws.Cells.Item(17, 1).Value2  = 4
$ws.Cells.Item(17, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value2  = "Los Lagos"
$ws.Cells.Item(17, 4).Value2  = 44749
$ws.Cells.Item(17, 5).Value2  = 10
$ws.Cells.Item(17, 6).Value2  = 100112052
$ws.Cells.Item(17, 7).Value2  = "Albahaca"
$ws.Cells.Item(17, 8).Value2  = "Sin especificar"
$ws.Cells.Item(17, 9).Value2  = "Primera"
$ws.Cells.Item(17, 10).Value2 = 50
$ws.Cells.Item(17, 11).Value2 = 7000
$ws.Cells.Item(17, 12).Value2 = 7000
$ws.Cells.Item(17, 13).Value2 = 7000
$ws.Cells.Item(17, 14).Value2 = "$/paquete"
$ws.Cells.Item(17, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value2 = 7000
$ws.Cells.Item(17, 17).Value2 = 1
$ws.Cells.Item(17, 18).Value2 = "Hortaliza"
